$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B5 holds the "TreatmentTab" query text. The query wrapped REPLACE() in a
# redundant CONCAT() call - strip that wrapper while leaving the rest of the
# query text (and all other cells) untouched.
$old = $ws.Range("B5").Value()
$new = $old.Replace("CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", "REPLACE(trt.treatment_agent, ';', ', ')")
$ws.Range("B5").Value = $new

# The saved view state moved up one row: the window's top-left cell and the
# active selection shifted from row 6 to row 5.
$ws.Range("C5").Select()
$excel.ActiveWindow.ScrollRow = 5
